$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 652 (shifts existing rows 652:766 down to 653:767)
$ws.Rows.Item(652).Insert()

# Populate the newly inserted row with the new weekly price observation
$ws.Cells.Item(652, 1).Value = 6
$ws.Cells.Item(652, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(652, 3).Value = "Metropolitana"
$ws.Cells.Item(652, 4).Value = 45180
$ws.Cells.Item(652, 5).Value = 13
$ws.Cells.Item(652, 6).Value = 100112039
$ws.Cells.Item(652, 7).Value = "Ciboulette"
$ws.Cells.Item(652, 8).Value = "Sin especificar"
$ws.Cells.Item(652, 9).Value = "Primera"
$ws.Cells.Item(652, 10).Value = 470
$ws.Cells.Item(652, 11).Value = 900
$ws.Cells.Item(652, 12).Value = 1000
$ws.Cells.Item(652, 13).Value = 953
$ws.Cells.Item(652, 14).Value = "`$/docena de atados"
$ws.Cells.Item(652, 15).Value = "Región Metropolitana"
$ws.Cells.Item(652, 16).Value = 318
$ws.Cells.Item(652, 17).Value = 3
$ws.Cells.Item(652, 18).Value = "Hortaliza"
